$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G3").Value = 500
$ws.Range("G4").Value = 1300
$ws.Range("G8").Value = 2500
$ws.Range("G9").Value = 1500
$ws.Range("G10").Value = 900
$ws.Range("G12").Value = 1500
$ws.Range("G15").Value = 1000
$ws.Range("G20").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 500
$ws.Range("G27").Value = 5000
$ws.Range("G30").Value = 0
$ws.Range("G31").Value = 1000
$ws.Range("G33").Value = 27750
